# Refreshes the hourly crypto snapshot table on Sheet1 (rows 2-51, one coin
# per row: B=Coin, C=Link, D=Price, E=Volume(1h)) to match the latest pull.
# Most rows just get new Price / Volume(1h) figures; rows 42-45 additionally
# rotate which coin occupies which row (RenderToken/Bittensor/Filecoin/Aave).
#
# Column D is persisted as text (t="inlineStr") even when the price reads like
# a plain number ("1.00", "5.26", ...). Writing such a string straight to
# .Value would let Excel auto-coerce it to a real number - silently dropping
# the trailing zero and introducing binary floating-point noise (1.00 -> 1,
# 5.26 -> 5.2599999999999998). Prefixing the value with a leading apostrophe
# is Excel's normal quote-prefix / "force text" convention: the cell keeps
# exactly the authored text, and the apostrophe itself never shows up in
# .Value / .Text or gets written to the file. Columns B, C and E are already
# unambiguous text (links, names, space-padded "  +x.xx%  " strings) so they
# are assigned directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '''58.830.50'
$ws.Range('E2').Value = '  +2.43%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '''2.521.82'
$ws.Range('E3').Value = '  +3.89%  '

# Row 4: TetherUSD
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  -0.01%  '

# Row 5: BNB
$ws.Range('D5').Value = '''534.84'
$ws.Range('E5').Value = '  +6.24%  '

# Row 6: Solana
$ws.Range('D6').Value = '''134.55'
$ws.Range('E6').Value = '  +5.34%  '

# Row 7: USDC
$ws.Range('D7').Value = '''0.999'
$ws.Range('E7').Value = '  +0.00%  '

# Row 8: XRP
$ws.Range('E8').Value = '  +3.72%  '

# Row 9: LidoStakedEther
$ws.Range('D9').Value = '''2.519.77'
$ws.Range('E9').Value = '  +3.39%  '

# Row 10: Dogecoin
$ws.Range('D10').Value = '''0.0996'
$ws.Range('E10').Value = '  +5.37%  '

# Row 11: TRON
$ws.Range('E11').Value = '  -1.27%  '

# Row 12: Toncoin
$ws.Range('D12').Value = '''5.26'
$ws.Range('E12').Value = '  +2.02%  '

# Row 13: Cardano
$ws.Range('D13').Value = '''0.333'
$ws.Range('E13').Value = '  +2.00%  '

# Row 14: WrappedliquidstakedEther2.0
$ws.Range('D14').Value = '''2.961.56'
$ws.Range('E14').Value = '  +3.56%  '

# Row 15: WrappedBTC
$ws.Range('D15').Value = '''58.777.67'
$ws.Range('E15').Value = '  +2.40%  '

# Row 16: Avalanche
$ws.Range('D16').Value = '''22.47'
$ws.Range('E16').Value = '  +3.95%  '

# Row 17: ShibaInu
$ws.Range('E17').Value = '  +3.98%  '

# Row 18: WrappedEther
$ws.Range('D18').Value = '''2.518.11'
$ws.Range('E18').Value = '  +3.27%  '

# Row 19: Chainlink
$ws.Range('D19').Value = '''10.67'
$ws.Range('E19').Value = '  +2.89%  '

# Row 20: Polkadot
$ws.Range('D20').Value = '''4.26'
$ws.Range('E20').Value = '  +4.67%  '

# Row 21: BitcoinCash
$ws.Range('D21').Value = '''321.61'
$ws.Range('E21').Value = '  +2.82%  '

# Row 22: Uniswap
$ws.Range('D22').Value = '''6.25'
$ws.Range('E22').Value = '  +10.33%  '

# Row 23: Dai
$ws.Range('D23').Value = '''1.00'
$ws.Range('E23').Value = '  +0.15%  '

# Row 24: Litecoin
$ws.Range('D24').Value = '''65.78'
$ws.Range('E24').Value = '  +4.36%  '

# Row 25: Polygon
$ws.Range('D25').Value = '''0.412'
$ws.Range('E25').Value = '  +2.41%  '

# Row 26: Binance-PegBSC-USD
$ws.Range('D26').Value = '''0.996'
$ws.Range('E26').Value = '  -0.43%  '

# Row 27: Kaspa
$ws.Range('E27').Value = '  +1.31%  '

# Row 28: InternetComputer(DFINITY)
$ws.Range('D28').Value = '''7.52'
$ws.Range('E28').Value = '  +5.39%  '

# Row 29: PEPE
$ws.Range('D29').Value = '''0.0₃0764'
$ws.Range('E29').Value = '  +6.73%  '

# Row 30: Monero
$ws.Range('D30').Value = '''172.44'
$ws.Range('E30').Value = '  +1.92%  '

# Row 31: PancakeSwap
$ws.Range('E31').Value = '  +6.15%  '

# Row 32: Fetch.AI
$ws.Range('E32').Value = '  +8.17%  '

# Row 33: Aptos
$ws.Range('D33').Value = '''6.36'
$ws.Range('E33').Value = '  +3.36%  '

# Row 34: USDe
$ws.Range('D34').Value = '''0.999'
$ws.Range('E34').Value = '  +0.04%  '

# Row 35: FirstDigitalUSD
$ws.Range('D35').Value = '''0.996'
$ws.Range('E35').Value = '  -0.09%  '

# Row 36: EthereumClassic
$ws.Range('D36').Value = '''18.20'
$ws.Range('E36').Value = '  +3.24%  '

# Row 37: ImmutableX
$ws.Range('E37').Value = '  +0.17%  '

# Row 38: NEARProtocol
$ws.Range('D38').Value = '''3.98'
$ws.Range('E38').Value = '  +2.91%  '

# Row 39: Stacks
$ws.Range('E39').Value = '  +5.77%  '

# Row 40: SuiNetwork
$ws.Range('D40').Value = '''0.822'
$ws.Range('E40').Value = '  +10.29%  '

# Row 41: OKB
$ws.Range('D41').Value = '''36.71'
$ws.Range('E41').Value = '  +0.78%  '

# Row 42: RenderToken -> Bittensor
$ws.Range('B42').Value = 'Bittensor'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D42').Value = '''278.00'
$ws.Range('E42').Value = '  +3.45%  '

# Row 43: Bittensor -> Filecoin
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').Value = '''3.49'
$ws.Range('E43').Value = '  +4.74%  '

# Row 44: Filecoin -> Aave
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').Value = '''131.84'
$ws.Range('E44').Value = '  +11.34%  '

# Row 45: Aave -> RenderToken
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').Value = '''5.07'
$ws.Range('E45').Value = '  +5.28%  '

# Row 46: Mantle
$ws.Range('D46').Value = '''0.594'
$ws.Range('E46').Value = '  +2.95%  '

# Row 47: Stellar
$ws.Range('E47').Value = '  +3.22%  '

# Row 48: Hedera
$ws.Range('D48').Value = '''0.0513'
$ws.Range('E48').Value = '  +6.41%  '

# Row 49: VeChain
$ws.Range('E49').Value = '  +6.17%  '

# Row 50: InjectiveProtocol
$ws.Range('D50').Value = '''17.11'
$ws.Range('E50').Value = '  +4.20%  '

# Row 51: Maker
$ws.Range('D51').Value = '''1.754.36'
$ws.Range('E51').Value = '  +3.67%  '
